# Update gh-pages output: refresh "want-to-go" counts (column F) and one
# venue address (D6) on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("D6").Value = "泗杨路张家港碧桂园天玺东南侧约60米 五月风华宴会中心"

$ws1.Range("F3").Value  = 571
$ws1.Range("F6").Value  = 1097
$ws1.Range("F7").Value  = 1440
$ws1.Range("F12").Value = 167
$ws1.Range("F14").Value = 437
$ws1.Range("F15").Value = 1358
$ws1.Range("F16").Value = 115
$ws1.Range("F17").Value = 111
$ws1.Range("F20").Value = 39
$ws1.Range("F21").Value = 655
$ws1.Range("F23").Value = 35
$ws1.Range("F24").Value = 225
$ws1.Range("F26").Value = 5886
$ws1.Range("F27").Value = 65
$ws1.Range("F31").Value = 14541
$ws1.Range("F32").Value = 1440
$ws1.Range("F33").Value = 213
$ws1.Range("F36").Value = 9071
$ws1.Range("F37").Value = 627
$ws1.Range("F38").Value = 4212

# --- Sheet "全部类型" (all types) ---
$ws2 = $wb.Worksheets.Item("全部类型")

$ws2.Range("D6").Value = "泗杨路张家港碧桂园天玺东南侧约60米 五月风华宴会中心"

$ws2.Range("F3").Value  = 571
$ws2.Range("F6").Value  = 1097
$ws2.Range("F7").Value  = 1440
$ws2.Range("F12").Value = 167
$ws2.Range("F14").Value = 437
$ws2.Range("F15").Value = 1358
$ws2.Range("F16").Value = 115
$ws2.Range("F17").Value = 111
$ws2.Range("F21").Value = 39
$ws2.Range("F22").Value = 655
$ws2.Range("F25").Value = 35
$ws2.Range("F26").Value = 225
$ws2.Range("F29").Value = 5886
$ws2.Range("F30").Value = 65
$ws2.Range("F34").Value = 14541
$ws2.Range("F35").Value = 1440
$ws2.Range("F36").Value = 213
$ws2.Range("F39").Value = 9071
$ws2.Range("F40").Value = 627
$ws2.Range("F41").Value = 4212
